$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.81"
$ws.Range("E2").Value = "'1.79%"

$ws.Range("D3").Value = "'37.55"
$ws.Range("E3").Value = "'0.17%"

$ws.Range("D4").Value = "'5.146"
$ws.Range("E4").Value = "'1.14%"

$ws.Range("D5").Value = "'0.07872"
$ws.Range("E5").Value = "'1.86%"

$ws.Range("D6").Value = "'4.420"
$ws.Range("E6").Value = "'1.45%"

$ws.Range("D7").Value = "'1.907"
$ws.Range("E7").Value = "'1.50%"

$ws.Range("D8").Value = "'8.275"
$ws.Range("E8").Value = "'0.84%"

$ws.Range("D9").Value = "'2.908"
$ws.Range("E9").Value = "'-8.98%"

$ws.Range("D10").Value = "'0.9181"
$ws.Range("E10").Value = "'0.08%"

$ws.Range("D11").Value = "'0.1188"
$ws.Range("E11").Value = "'0.37%"

$ws.Range("D12").Value = "'0.1915"
$ws.Range("E12").Value = "'1.84%"

$ws.Range("D13").Value = "'0.09076"
$ws.Range("E13").Value = "'4.03%"

$ws.Range("D14").Value = "'0.03348"
$ws.Range("E14").Value = "'-1.65%"

$ws.Range("D15").Value = "'0.09625"
$ws.Range("E15").Value = "'-0.94%"

$ws.Range("D16").Value = "'0.001396"
$ws.Range("E16").Value = "'2.31%"

$ws.Range("D17").Value = "'0.005716"
$ws.Range("E17").Value = "'-2.98%"

$ws.Range("D18").Value = "'3.517"
$ws.Range("E18").Value = "'-1.98%"

$ws.Range("D19").Value = "'0.3443"
$ws.Range("E19").Value = "'1.05%"

$ws.Range("D20").Value = "'5.274"
$ws.Range("E20").Value = "'5.14%"

$ws.Range("E21").Value = "'-0.24%"

$ws.Range("D22").Value = "'0.2593"
$ws.Range("E22").Value = "'-0.10%"

$ws.Range("D23").Value = "'0.04375"
$ws.Range("E23").Value = "'1.17%"

$ws.Range("D24").Value = "'0.001252"
$ws.Range("E24").Value = "'3.17%"

$ws.Range("D25").Value = "'0.004692"
$ws.Range("E25").Value = "'3.20%"

$ws.Range("D26").Value = "'0.0001366"
$ws.Range("E26").Value = "'1.02%"

$ws.Range("D27").Value = "'0.0003999"
$ws.Range("E27").Value = "'-98.10%"

$ws.Range("D39").Value = "'0.02297"
$ws.Range("E39").Value = "'3.43%"

$ws.Range("D40").Value = "'0.05074"
$ws.Range("E40").Value = "'3.12%"

$ws.Range("D41").Value = "'0.007478"
$ws.Range("E41").Value = "'-1.24%"

$ws.Range("D42").Value = "'0.009056"
$ws.Range("E42").Value = "'-7.90%"

$ws.Range("D43").Value = "'0.1351"
$ws.Range("E43").Value = "'1.05%"

$ws.Range("D44").Value = "'0.001959"
$ws.Range("E44").Value = "'-1.81%"

$ws.Range("D45").Value = "'0.008661"
$ws.Range("E45").Value = "'-1.65%"

$ws.Range("D46").Value = "'0.00006635"
$ws.Range("E46").Value = "'1.46%"

$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.27%"

$ws.Range("D48").Value = "'0.003361"
$ws.Range("E48").Value = "'12.07%"

$ws.Range("D49").Value = "'0.001002"
$ws.Range("E49").Value = "'-23.02%"

$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.27%"

$ws.Range("D51").Value = "'0.0002008"
$ws.Range("E51").Value = "'0.27%"
